# Generate Report for Handback
#
# The handback automation re-ran: the previously handed-back file
# (2d0bc8b8-8423-4709-9e31-1b4533162f81) was re-processed under a new
# GUID (c77dc451-783a-4ffe-b319-5fcc4e27c222) with refreshed timestamps,
# and a brand-new file (f2e2ca04-4b2a-42ce-9d37-9a0c66c7aa2c) was handed
# back as well, so every sheet grows from one data row to two.

$wb = $excel.ActiveWorkbook

$newGuid1 = "c77dc451-783a-4ffe-b319-5fcc4e27c222"
$newGuid2 = "f2e2ca04-4b2a-42ce-9d37-9a0c66c7aa2c"

$hash1 = "5826b1e6ee54820443aa4623f9ef8bfd61659d5b"
$hash2 = "38b128ea99c2c4b1467332e76b301841b4a5fc96"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2: refresh the GUID and the generate-date timestamp for the existing file
$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G2").Value = "2016-08-28 15:01:17"
$wsOverview.Range("G2").NumberFormat = $dateFmt

# Row 3: brand new handed-back file
$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-28 15:01:17"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25f63ee5604194b54bd9980102ef75139d21b599/e2e/$newGuid1.md", $null, $null, "e2e\$newGuid1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25f63ee5604194b54bd9980102ef75139d21b599/e2e/$newGuid2.md", $null, $null, "e2e\$newGuid2.md") | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 2: refresh GUID / hash / timestamps for the existing file
$wsZhCn.Range("A2").Value = "$newGuid1.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D2").Value = "e2e"
$wsZhCn.Range("E2").Value = "ht"
$wsZhCn.Range("F2").Value = "False"
$wsZhCn.Range("G2").Value = "$newGuid1.$hash1.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-28 15:01:13"
$wsZhCn.Range("H2").NumberFormat = $dateFmt
$wsZhCn.Range("I2").Value = "$newGuid1.md"
$wsZhCn.Range("J2").Value = "$newGuid1.$hash1.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-28 15:01:34"
$wsZhCn.Range("K2").NumberFormat = $dateFmt
$wsZhCn.Range("L2").Value = ""
$wsZhCn.Range("M2").Value = "True"
$wsZhCn.Range("N2").Value = ""
$wsZhCn.Range("O2").Value = "False"
$wsZhCn.Range("P2").Value = ""

# Row 3: brand new handed-back file
$wsZhCn.Range("A3").Value = "$newGuid2.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "$newGuid2.$hash2.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-28 15:01:13"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("I3").Value = "$newGuid2.md"
$wsZhCn.Range("J3").Value = "$newGuid2.$hash2.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-28 15:01:34"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25f63ee5604194b54bd9980102ef75139d21b599/e2e/$newGuid1.md", $null, $null, "$newGuid1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/682f20195f1ad936bf36f60e462a365444dba15c/e2e/$newGuid1.md", $null, $null, "$newGuid1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25f63ee5604194b54bd9980102ef75139d21b599/e2e/$newGuid2.md", $null, $null, "$newGuid2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/682f20195f1ad936bf36f60e462a365444dba15c/e2e/$newGuid2.md", $null, $null, "$newGuid2.md") | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 2: refresh GUID / hash / timestamps for the existing file
$wsDeDe.Range("A2").Value = "$newGuid1.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D2").Value = "e2e"
$wsDeDe.Range("E2").Value = "ht"
$wsDeDe.Range("F2").Value = "False"
$wsDeDe.Range("G2").Value = "$newGuid1.$hash1.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-28 15:01:17"
$wsDeDe.Range("H2").NumberFormat = $dateFmt
$wsDeDe.Range("I2").Value = "$newGuid1.md"
$wsDeDe.Range("J2").Value = "$newGuid1.$hash1.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-28 15:01:41"
$wsDeDe.Range("K2").NumberFormat = $dateFmt
$wsDeDe.Range("L2").Value = ""
$wsDeDe.Range("M2").Value = "True"
$wsDeDe.Range("N2").Value = ""
$wsDeDe.Range("O2").Value = "False"
$wsDeDe.Range("P2").Value = ""

# Row 3: brand new handed-back file
$wsDeDe.Range("A3").Value = "$newGuid2.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "$newGuid2.$hash2.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-28 15:01:17"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("I3").Value = "$newGuid2.md"
$wsDeDe.Range("J3").Value = "$newGuid2.$hash2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-28 15:01:41"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25f63ee5604194b54bd9980102ef75139d21b599/e2e/$newGuid1.md", $null, $null, "$newGuid1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bffd5d1046dcebe704b249925e52d2bb93585149/e2e/$newGuid1.md", $null, $null, "$newGuid1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25f63ee5604194b54bd9980102ef75139d21b599/e2e/$newGuid2.md", $null, $null, "$newGuid2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bffd5d1046dcebe704b249925e52d2bb93585149/e2e/$newGuid2.md", $null, $null, "$newGuid2.md") | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
